$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Players")
$ws2 = $wb.Worksheets.Item("OwnerTotals")

# Column G width change: 35 -> 18
# (the engine's width round-trip adds ~0.8333 to whatever ColumnWidth is set,
# so we back that out here to land on exactly 18 in the saved OOXML)
$ws.Columns.Item(7).ColumnWidth = 17.166666666666668

# Players sheet cell updates
$ws.Cells.Item(2, 7).Value = "15:41 - 2nd Half"
$ws.Cells.Item(2, 10).Value = 6
$ws.Cells.Item(2, 15).Value = 2
$ws.Cells.Item(2, 16).Value = 17
$ws.Cells.Item(2, 18).Value = 8
$ws.Cells.Item(2, 20).Value = 3
$ws.Cells.Item(3, 7).Value = "15:41 - 2nd Half"
$ws.Cells.Item(3, 8).Value = 15
$ws.Cells.Item(3, 9).Value = 16
$ws.Cells.Item(3, 16).Value = 23
$ws.Cells.Item(3, 17).Value = 7
$ws.Cells.Item(3, 18).Value = 12
$ws.Cells.Item(3, 20).Value = 4
$ws.Cells.Item(4, 7).Value = "15:41 - 2nd Half"
$ws.Cells.Item(4, 8).Value = 19
$ws.Cells.Item(4, 9).Value = 20
$ws.Cells.Item(4, 10).Value = 2
$ws.Cells.Item(4, 13).Value = 1
$ws.Cells.Item(4, 16).Value = 24
$ws.Cells.Item(4, 17).Value = 7
$ws.Cells.Item(4, 18).Value = 11
$ws.Cells.Item(4, 19).Value = 5
$ws.Cells.Item(4, 20).Value = 9
$ws.Cells.Item(5, 7).Value = "15:41 - 2nd Half"
$ws.Cells.Item(5, 8).Value = 1
$ws.Cells.Item(5, 9).Value = 2
$ws.Cells.Item(5, 10).Value = 2
$ws.Cells.Item(5, 16).Value = 15
$ws.Cells.Item(5, 17).Value = 1
$ws.Cells.Item(5, 18).Value = 5
$ws.Cells.Item(5, 20).Value = 4
$ws.Cells.Item(6, 4).Value = "Tyler Tanner"
$ws.Cells.Item(6, 5).Value = "VAN"
$ws.Cells.Item(6, 7).Value = "15:41 - 2nd Half"
$ws.Cells.Item(6, 8).Value = 6
$ws.Cells.Item(6, 10).Value = 3
$ws.Cells.Item(6, 11).Value = 4
$ws.Cells.Item(6, 12).Value = 3
$ws.Cells.Item(6, 15).Value = 3
$ws.Cells.Item(6, 16).Value = 26
$ws.Cells.Item(6, 18).Value = 7
$ws.Cells.Item(6, 20).Value = 2
$ws.Cells.Item(7, 4).Value = "Rashaun Agee"
$ws.Cells.Item(7, 5).Value = "TA&M"
$ws.Cells.Item(7, 7).Value = "15:41 - 2nd Half"
$ws.Cells.Item(7, 8).Value = 5
$ws.Cells.Item(7, 9).Value = 3
$ws.Cells.Item(7, 10).Value = 6
$ws.Cells.Item(7, 11).Value = 0
$ws.Cells.Item(7, 12).Value = 0
$ws.Cells.Item(7, 13).Value = 2
$ws.Cells.Item(7, 16).Value = 18
$ws.Cells.Item(7, 20).Value = 0
$ws.Cells.Item(7, 21).Value = 1
$ws.Cells.Item(7, 22).Value = 2
$ws.Cells.Item(8, 7).Value = "15:41 - 2nd Half"
$ws.Cells.Item(8, 10).Value = 7
$ws.Cells.Item(8, 14).Value = 2
$ws.Cells.Item(8, 16).Value = 17
$ws.Cells.Item(8, 18).Value = 2
$ws.Cells.Item(8, 20).Value = 1
$ws.Cells.Item(9, 7).Value = "15:41 - 2nd Half"
$ws.Cells.Item(9, 8).Value = -1
$ws.Cells.Item(9, 10).Value = 1
$ws.Cells.Item(9, 15).Value = 1
$ws.Cells.Item(9, 16).Value = 17
$ws.Cells.Item(9, 18).Value = 5
$ws.Cells.Item(9, 20).Value = 4
$ws.Cells.Item(10, 7).Value = "15:41 - 2nd Half"
$ws.Cells.Item(10, 8).Value = 15
$ws.Cells.Item(10, 9).Value = 12
$ws.Cells.Item(10, 10).Value = 3
$ws.Cells.Item(10, 11).Value = 3
$ws.Cells.Item(10, 14).Value = 1
$ws.Cells.Item(10, 16).Value = 22
$ws.Cells.Item(10, 17).Value = 4
$ws.Cells.Item(10, 18).Value = 7
$ws.Cells.Item(11, 7).Value = "15:41 - 2nd Half"
$ws.Cells.Item(11, 8).Value = 9
$ws.Cells.Item(11, 10).Value = 6
$ws.Cells.Item(11, 12).Value = 1
$ws.Cells.Item(11, 16).Value = 13
$ws.Cells.Item(11, 18).Value = 7
$ws.Cells.Item(12, 4).Value = "Chandler Bing"
$ws.Cells.Item(12, 5).Value = "VAN"
$ws.Cells.Item(12, 7).Value = "15:41 - 2nd Half"
$ws.Cells.Item(12, 8).Value = 9
$ws.Cells.Item(12, 9).Value = 7
$ws.Cells.Item(12, 10).Value = 5
$ws.Cells.Item(12, 13).Value = 1
$ws.Cells.Item(12, 14).Value = 0
$ws.Cells.Item(12, 15).Value = 1
$ws.Cells.Item(12, 16).Value = 20
$ws.Cells.Item(12, 17).Value = 3
$ws.Cells.Item(12, 18).Value = 7
$ws.Cells.Item(12, 19).Value = 1
$ws.Cells.Item(12, 20).Value = 3
$ws.Cells.Item(13, 4).Value = "Federiko Federiko"
$ws.Cells.Item(13, 5).Value = "TA&M"
$ws.Cells.Item(13, 7).Value = "15:41 - 2nd Half"
$ws.Cells.Item(13, 8).Value = 5
$ws.Cells.Item(13, 9).Value = 0
$ws.Cells.Item(13, 10).Value = 4
$ws.Cells.Item(13, 12).Value = 1
$ws.Cells.Item(13, 16).Value = 8
$ws.Cells.Item(13, 17).Value = 0
$ws.Cells.Item(13, 18).Value = 0
$ws.Cells.Item(13, 20).Value = 0
$ws.Cells.Item(14, 4).Value = "Zach Clemence"
$ws.Cells.Item(14, 7).Value = "15:41 - 2nd Half"
$ws.Cells.Item(14, 8).Value = 5
$ws.Cells.Item(14, 9).Value = 6
$ws.Cells.Item(14, 14).Value = 2
$ws.Cells.Item(14, 15).Value = 4
$ws.Cells.Item(14, 16).Value = 12
$ws.Cells.Item(14, 17).Value = 3
$ws.Cells.Item(14, 18).Value = 4
$ws.Cells.Item(14, 20).Value = 1
$ws.Cells.Item(15, 7).Value = "15:41 - 2nd Half"
$ws.Cells.Item(15, 8).Value = 2
$ws.Cells.Item(15, 11).Value = 1
$ws.Cells.Item(15, 15).Value = 4
$ws.Cells.Item(15, 16).Value = 13
$ws.Cells.Item(16, 7).Value = "15:41 - 2nd Half"
$ws.Cells.Item(16, 11).Value = 1
$ws.Cells.Item(16, 15).Value = 1
$ws.Cells.Item(16, 16).Value = 17
$ws.Cells.Item(16, 18).Value = 5
$ws.Cells.Item(16, 20).Value = 3
$ws.Cells.Item(17, 7).Value = "15:41 - 2nd Half"
$ws.Cells.Item(18, 7).Value = "15:41 - 2nd Half"
$ws.Cells.Item(19, 7).Value = "15:41 - 2nd Half"
$ws.Cells.Item(19, 8).Value = -3
$ws.Cells.Item(19, 16).Value = 14
$ws.Cells.Item(19, 18).Value = 4
$ws.Cells.Item(20, 7).Value = "15:41 - 2nd Half"
$ws.Cells.Item(20, 16).Value = 4

# OwnerTotals sheet cell updates
$ws2.Cells.Item(2, 1).Value = "Hilton Heads"
$ws2.Cells.Item(2, 2).Value = 19
$ws2.Cells.Item(3, 1).Value = "Boozers Losers"
$ws2.Cells.Item(3, 2).Value = 13
$ws2.Cells.Item(4, 2).Value = 11
$ws2.Cells.Item(5, 2).Value = 1
